$d = $word.ActiveDocument

$d.Content.Find.Execute("670×8=5360", $true, $false, $false, $false, $false, $true, 1, $false, "217×7=1519", 2)
$d.Content.Find.Execute("206×5=1030", $true, $false, $false, $false, $false, $true, 1, $false, "682×8=5456", 2)
$d.Content.Find.Execute("427×7=2989", $true, $false, $false, $false, $false, $true, 1, $false, "239×8=1912", 2)
$d.Content.Find.Execute("847×2=1694", $true, $false, $false, $false, $false, $true, 1, $false, "654×5=3270", 2)
$d.Content.Find.Execute("953×3=2859", $true, $false, $false, $false, $false, $true, 1, $false, "460×7=3220", 2)
$d.Content.Find.Execute("519×5=2595", $true, $false, $false, $false, $false, $true, 1, $false, "279×4=1116", 2)
$d.Content.Find.Execute("265×9=2385", $true, $false, $false, $false, $false, $true, 1, $false, "210×6=1260", 2)
$d.Content.Find.Execute("843×9=7587", $true, $false, $false, $false, $false, $true, 1, $false, "384×3=1152", 2)
$d.Content.Find.Execute("155×3=465", $true, $false, $false, $false, $false, $true, 1, $false, "367×6=2202", 2)
$d.Content.Find.Execute("365×3=1095", $true, $false, $false, $false, $false, $true, 1, $false, "366×4=1464", 2)
$d.Content.Find.Execute("657×6=3942", $true, $false, $false, $false, $false, $true, 1, $false, "157×7=1099", 2)
$d.Content.Find.Execute("626×4=2504", $true, $false, $false, $false, $false, $true, 1, $false, "674×6=4044", 2)
$d.Content.Find.Execute("931×6=5586", $true, $false, $false, $false, $false, $true, 1, $false, "386×9=3474", 2)
$d.Content.Find.Execute("970×6=5820", $true, $false, $false, $false, $false, $true, 1, $false, "351×5=1755", 2)
$d.Content.Find.Execute("585×3=1755", $true, $false, $false, $false, $false, $true, 1, $false, "761×3=2283", 2)
$d.Content.Find.Execute("113×9=1017", $true, $false, $false, $false, $false, $true, 1, $false, "740×6=4440", 2)
$d.Content.Find.Execute("443×6=2658", $true, $false, $false, $false, $false, $true, 1, $false, "159×6=954", 2)
$d.Content.Find.Execute("321×6=1926", $true, $false, $false, $false, $false, $true, 1, $false, "842×5=4210", 2)
$d.Content.Find.Execute("187×3=561", $true, $false, $false, $false, $false, $true, 1, $false, "818×3=2454", 2)
$d.Content.Find.Execute("799×7=5593", $true, $false, $false, $false, $false, $true, 1, $false, "968×5=4840", 2)
$d.Content.Find.Execute("198×8=1584", $true, $false, $false, $false, $false, $true, 1, $false, "879×2=1758", 2)
$d.Content.Find.Execute("314×8=2512", $true, $false, $false, $false, $false, $true, 1, $false, "466×2=932", 2)
$d.Content.Find.Execute("971×5=4855", $true, $false, $false, $false, $false, $true, 1, $false, "214×6=1284", 2)
$d.Content.Find.Execute("588×5=2940", $true, $false, $false, $false, $false, $true, 1, $false, "180×3=540", 2)
$d.Content.Find.Execute("866×6=5196", $true, $false, $false, $false, $false, $true, 1, $false, "205×7=1435", 2)
